$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header: Active Terminals count 22 -> 23
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = 23

# ---------------------------------------------------------------------------
# 2. Insert a new row at 27 (old row 27 "Total Outstanding..." shifts to 28,
#    and merged cells / dimension follow automatically).
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Insert()

# ---------------------------------------------------------------------------
# 3. Helper to fully (re)write one terminal data row (columns A,C,E,F,H,I,J,K,L,M,N,O,P,Q)
#    hasH / hasI control whether the optional Est.Cash Out / Last Error cells exist.
#    (Positional parameters only -- named parameter binding isn't reliable here.)
# ---------------------------------------------------------------------------
function Set-TermRow($r, $A, $C, $E, $F, $hasH, $H, $hasI, $I, $J, $K, $L, $M, $N, $O, $P, $Q) {
    $ws.Cells.Item($r, 1).Value = $A
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    if ($hasH) {
        $ws.Cells.Item($r, 8).Value = $H
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
    if ($hasI) {
        $ws.Cells.Item($r, 9).Value = $I
    } else {
        $ws.Cells.Item($r, 9).ClearContents()
    }
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 13).Value = $M
    $ws.Cells.Item($r, 14).Value = $N
    $ws.Cells.Item($r, 15).Value = $O
    $ws.Cells.Item($r, 16).Value = $P
    $ws.Cells.Item($r, 17).Value = $Q
}

# ---------------------------------------------------------------------------
# 4. Rewrite every terminal row (5-27) with the post-edit values.
# ---------------------------------------------------------------------------
Set-TermRow 5  "LK644532" "SCL ENTERPRISES LAUNDRY" 1880  "T" $true  45127.042061423606 $false "" "06/28/23 15:33" "06/28/23 15:33" 0   "`$1,880 as of 6/28/2023 1:33:05 PM"  1980  0 0 0
Set-TermRow 6  "L647934"  "SB #6"                   1940  "T" $false 0                  $true  "ATM Inactive greater than 2000 minutes" "04/06/23 22:10" "04/06/23 22:05" 20 "`$1,940 as of 4/6/2023 8:05:45 PM"  1960  0 0 0
Set-TermRow 7  "L662336"  "SB#4 MONA MARKET"        2240  "T" $true  45114.042061423606 $true  "The Triton ATM reported error code 35: Mistracked Note at Double Detec" "06/29/23 15:48" "06/29/23 15:48" 0 "`$2,300 as of 6/29/2023 8:00:28 AM" 2240 0 0 0
Set-TermRow 8  "L474746"  "ZACATES MARKET"          2360  "T" $true  45114.042061423606 $false "" "06/29/23 19:45" "06/29/23 19:45" 0 "`$2,480 as of 6/29/2023 9:34:21 AM" 2480 0 0 0
Set-TermRow 9  "L682801"  "SB#5"                    2520  "T" $true  45155.042061423606 $false "" "06/28/23 18:23" "06/28/23 18:23" 0 "`$2,520 as of 6/28/2023 4:23:23 PM" 2560 0 0 0
Set-TermRow 10 "LK561655" "CRENSHAW CRAVOR #2"      2780  "T" $false 0                  $true  "ATM Inactive greater than 48 minutes" "01/23/20 08:24" "01/23/20 08:24" 0 "`$2,780 as of 1/23/2020 6:24:32 AM" 2800 0 0 0
Set-TermRow 11 "L474792"  "NICK SHELL SERVICE"      3440  "T" $true  45142.042061423606 $false "" "06/29/23 17:17" "06/29/23 17:17" 0 "`$3,500 as of 6/28/2023 5:44:16 PM" 3460 0 0 0
Set-TermRow 12 "L475182"  "LA ESQUINA DE ORO"       3800  "T" $false 0                  $true  "ATM Inactive greater than 48 minutes" "09/16/20 16:57" "09/15/20 23:38" 0 "`$3,800 as of 9/16/2020 1:28:00 PM" 3800 0 0 0
Set-TermRow 13 "L688961"  "MONA MART"               3900  "T" $true  45140.042061423606 $false "" "06/29/23 16:20" "06/29/23 16:20" 0 "`$1,940 as of 6/28/2023 4:57:51 PM" 3940 0 0 0
Set-TermRow 14 "L475090"  "S.B. 2"                  4680  "T" $true  45115.042061423606 $false "" "06/29/23 19:04" "06/29/23 19:04" 0 "`$4,900 as of 6/29/2023 11:52:24 AM" 4720 0 0 0
Set-TermRow 15 "LK923383" "SAMYS PHONE CARDS"       4860  "T" $true  45112.042061423606 $false "" "06/29/23 18:34" "06/29/23 18:34" 100 "`$5,680 as of 6/28/2023 8:26:26 PM" 4900 0 0 0
Set-TermRow 16 "L476340"  "DONUT & SANDWICH"        4940  "T" $true  45139.042061423606 $false "" "06/29/23 16:45" "06/29/23 16:45" 0 "`$4,300 as of 6/24/2023 11:23:13 AM" 4960 0 0 0
Set-TermRow 17 "LK236828" "WORLDWIDE AUTOMOTIVE"    4940  "T" $true  45129.042061423606 $false "" "06/29/23 12:32" "06/29/23 11:23" 80 "`$4,940 as of 6/29/2023 10:32:07 AM" 5140 0 0 0
Set-TermRow 18 "L474761"  "BABS MARKET"             5160  "T" $true  45168.042061423606 $false "" "06/29/23 19:19" "06/29/23 19:19" 40 "`$5,300 as of 6/29/2023 10:52:54 AM" 5260 0 0 0
Set-TermRow 19 "L697589"  "S B DISCOUNT MART"       5780  "T" $true  45114.042061423606 $false "" "06/29/23 19:31" "06/29/23 19:31" 20 "`$7,180 as of 6/29/2023 11:54:34 AM" 5880 0 0 0
Set-TermRow 20 "LK864765" "SKY LIQUOR"              5940  "T" $true  45121.042061423606 $false "" "06/29/23 19:02" "06/29/23 11:48" 0 "`$5,960 as of 6/29/2023 9:48:28 AM" 5940 0 0 0
Set-TermRow 21 "L488595"  "N S MART"                6020  "T" $true  45407.042061423606 $false "" "06/28/23 19:23" "06/28/23 19:23" 0 "`$6,020 as of 6/28/2023 5:23:06 PM" 6180 0 0 0
Set-TermRow 22 "L474817"  "SAFETY MARKET"           6440  "T" $true  45123.042061423606 $false "" "06/29/23 17:44" "06/29/23 17:44" 100 "`$6,480 as of 6/29/2023 7:54:04 AM" 6460 0 0 0
Set-TermRow 23 "L688966"  "LACON MINI MART"         6900  "T" $true  45279.042061423606 $true  "ATM Inactive greater than 2000 minutes" "06/25/23 11:10" "06/22/23 16:56" 20 "`$6,900 as of 6/22/2023 2:56:56 PM" 6920 0 0 0
Set-TermRow 24 "L704741"  "W ADAMS COIN LAUNDRY"    8320  "T" $false 0                  $false "" "06/29/23 19:25" "06/29/23 19:25" 0 "`$8,680 as of 6/29/2023 11:29:36 AM" 8460 0 0 0
Set-TermRow 25 "L697590"  "S B MARKET ST"           8780  "T" $true  45120.042061423606 $false "" "06/29/23 11:36" "06/29/23 11:36" 0 "`$8,780 as of 6/29/2023 9:36:36 AM" 8800 0 0 0
Set-TermRow 26 "L678988"  "PAYELESS MARKET"         9040  "T" $true  45140.042061423606 $false "" "06/29/23 17:14" "06/29/23 14:45" 0 "`$9,240 as of 6/29/2023 10:47:16 AM" 9040 0 0 0
Set-TermRow 27 "LK891176" "98 DISCOUNT STORE"       20540 "T" $true  45122.042061423606 $false "" "06/29/23 18:49" "06/29/23 18:49" 0 "`$21,100 as of 6/29/2023 10:23:59 AM" 20580 0 0 0

# ---------------------------------------------------------------------------
# 5. Total row, now at row 28.
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Total Outstanding Cash Balance:"
$ws.Range("E28").Value = 127200

# ---------------------------------------------------------------------------
# 6. Merge cells for the new row 27 (matching the A:B / C:D / F:G pattern used
#    by every other data row), then re-merge the total row so the merge list
#    ends up in the same order as the target (row 27 merges, then row 28).
# ---------------------------------------------------------------------------
$ws.Range("A28:D28").UnMerge()
$ws.Range("E28:H28").UnMerge()
$ws.Range("A27:B27").Merge()
$ws.Range("C27:D27").Merge()
$ws.Range("F27:G27").Merge()
$ws.Range("A28:D28").Merge()
$ws.Range("E28:H28").Merge()

# ---------------------------------------------------------------------------
# 7. Style metadata: numFmt 166 / 167 formatCode swap (cosmetic, unused by any
#    cell in this sheet -- applied best effort via the object model).
# ---------------------------------------------------------------------------
$wb.DeleteNumberFormat("[$-010409]`$#,##0")
$wb.DeleteNumberFormat("[$-010409]m/d/yyyy")
